# NIT-8110388328.xlsx edit: remove previous "Estado de Cuenta" worker rows,
# add the new worker mora records, update the summary header counts/total,
# and push the signature block down to make room for the expanded table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Make room: the worker table used to be a single data row (16);
#    it grows to 7 rows (16-22). Insert 6 rows below row 16 so the
#    trailing "firma" / signature block (previously rows 21-22) is
#    pushed down to rows 27-28, and merged cells shift with it.
# ------------------------------------------------------------------
$ws.Rows("17:22").Insert()

# ------------------------------------------------------------------
# 2) Give the newly inserted rows the same look as the existing data
#    row (borders/fonts/number formats), then fill in the values.
# ------------------------------------------------------------------
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 16 - JOSE ALEJANDRO MARANTO ARROYO, periodo 2506
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1048609587"
$ws.Range("D16").Value = "JOSE ALEJANDRO MARANTO ARROYO"
$ws.Range("E16").Value = "2506"
$ws.Range("F16").Value = 28470
$ws.Range("G16").Value = 1423500

# Row 17 - JHON JAIRO PEREZ PAJARO, periodo 2507
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1001899444"
$ws.Range("D17").Value = "JHON JAIRO PEREZ PAJARO"
$ws.Range("E17").Value = "2507"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# Row 18 - JHON JAIRO PEREZ PAJARO, periodo 2506
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1001899444"
$ws.Range("D18").Value = "JHON JAIRO PEREZ PAJARO"
$ws.Range("E18").Value = "2506"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

# Row 19 - CESAR MARRUGO GUERRERO, periodo 2507
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "9185622"
$ws.Range("D19").Value = "CESAR MARRUGO GUERRERO"
$ws.Range("E19").Value = "2507"
$ws.Range("F19").Value = 72000
$ws.Range("G19").Value = 1800000

# Row 20 - CESAR MARRUGO GUERRERO, periodo 2506
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "9185622"
$ws.Range("D20").Value = "CESAR MARRUGO GUERRERO"
$ws.Range("E20").Value = "2506"
$ws.Range("F20").Value = 72000
$ws.Range("G20").Value = 1800000

# Row 21 - JOSE DANIEL ALVAREZ TOVAR, periodo 2506 (previously the lone row 16)
$ws.Range("B21").Value = "PE"
$ws.Range("C21").Value = "726342302121993"
$ws.Range("D21").Value = "JOSE DANIEL ALVAREZ TOVAR"
$ws.Range("E21").Value = "2506"
$ws.Range("F21").Value = 36341
$ws.Range("G21").Value = 908526

# Row 22 - ESNEIDER SEGUNDO GOMEZ ALEAN, periodo 2506
$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1003589833"
$ws.Range("D22").Value = "ESNEIDER SEGUNDO GOMEZ ALEAN"
$ws.Range("E22").Value = "2506"
$ws.Range("F22").Value = 22776
$ws.Range("G22").Value = 1423500

# ------------------------------------------------------------------
# 3) Update the header summary block: total mora, worker count and
#    period count.
# ------------------------------------------------------------------
$ws.Range("E11").Value = 345467
$ws.Range("C13").Value = 5
$ws.Range("F13").Value = 2

# Column D ("Nombre Trabajador") needs to accommodate the longest new
# name ("JOSE ALEJANDRO MARANTO ARROYO"), so widen it to fit.
$ws.Columns("D:D").ColumnWidth = 34.1

Write-Host "edit applied"
